$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$B2Text = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
    prg.website AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@
$ws.Range("B2").Value = $B2Text

$C2Text = @'
SELECT DISTINCT
    COUNT(DISTINCT prg.program_id) AS "Programs",
    COUNT(DISTINCT prj.project_id) AS "Projects",
    COUNT(DISTINCT gnt.grant_id) AS "Grants",
    COUNT(DISTINCT pub.pmid) AS "Publications"
FROM 
    df_program prg
LEFT JOIN 
    df_project prj ON prg.program_id = prj."program.program_id"
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Cervical Cancer%';
'@
$ws.Range("C2").Value = $C2Text

$B3Text = @'
SELECT DISTINCT
    prj.project_id AS "Project ID", 
    prj.project_title AS "Project Title",
    prj.org_name AS "Organization",
    prj.project_start_date AS "Project Start Date",
    prj.project_end_date AS "Project End Date"
FROM 
    df_project prj
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
   lower(prj.project_id) ASC
LIMIT 100;
'@
$ws.Range("B3").Value = $B3Text

$B4Text = @'
SELECT DISTINCT
    gnt.grant_id AS "Grant ID", 
    prj.project_id AS "Project",
    gnt.grant_title AS "Grant Title",
    gnt.principal_investigators AS "Principal Investigators",
    gnt.program_officers AS "Program Officers",
    gnt.fiscal_year AS "Fiscal Year",
    gnt.project_end_date AS "Project End Date"
FROM 
    df_grant gnt
LEFT JOIN 
    df_project prj ON gnt."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_publication pub ON prj.project_id = pub."project.project_id"
WHERE 
    prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    lower(gnt.grant_id) ASC
LIMIT 100;
'@
$ws.Range("B4").Value = $B4Text

$B5Text = @'
SELECT DISTINCT
    pub.pmid AS "PubMed ID", 
    pub.title AS "Title",
    pub.authors AS "Authors",
    pub.publication_date AS "Publication Date",
    pub.cited_by AS "Cited By",
    CASE 
    WHEN pub.relative_citation_ratio = 0 THEN '0'
    WHEN pub.relative_citation_ratio = 7.0 THEN '7'
    WHEN pub.relative_citation_ratio = 2.0 THEN '2'
  WHEN pub.relative_citation_ratio = 3.0 THEN '3'
    WHEN pub.relative_citation_ratio = ROUND(pub.relative_citation_ratio) THEN CAST(ROUND(pub.relative_citation_ratio) AS VARCHAR) 
    ELSE CAST(ROUND(pub.relative_citation_ratio, 2) AS VARCHAR)
END AS "Relative Citation Ratio"
FROM 
    df_publication pub
LEFT JOIN 
    df_project prj ON pub."project.project_id" = prj.project_id
LEFT JOIN 
    df_program prg ON prj."program.program_id" = prg.program_id
LEFT JOIN 
    df_grant gnt ON prj.project_id = gnt."project.project_id"
WHERE 
     prg.cancer_type LIKE '%Cervical Cancer%'
ORDER BY 
    pub.pmid ASC
LIMIT 100;
'@
$ws.Range("B5").Value = $B5Text

# Update the active selection to C5, matching the final user interaction
$ws.Range("C5").Select()
